$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.825.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.80%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.572.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.07%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = "'563.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.27%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'142.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.23%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.73%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.575.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.11%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'6.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.15%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +3.68%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +7.95%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +3.23%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.024.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.11%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'58.911.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.04%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'21.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +7.31%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +5.59%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.580.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.66%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.50%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'335.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.67%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.94%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +2.07%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.00%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'63.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.34%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.447"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +6.52%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.07%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +2.46%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +2.94%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0₃0781"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +6.90%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.00%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +1.76%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'159.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.87%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'6.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.45%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'18.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.93%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +3.13%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +4.25%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +8.20%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +3.90%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'36.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.73%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +5.05%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'291.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +5.12%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +2.12%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.05%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +3.26%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.75%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -0.27%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0532"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.58%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'19.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.44%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'124.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +12.36%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0231"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.25%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +5.37%  "
$ws.Range("E51").Style = "Normal"
